{"js": "// Commit: \"adding in appendix for data prep tools\"\n//\n// Appends a new \"Appendix\" heading (style H1) followed by an empty\n// paragraph (style \"Para\") to the very end of the document body, after\n// the last existing paragraph (\"...Reiterating that data analysis\n// assists our thinking, it does not replace it\").\n//\n// The appendix heading text is built from three logical pieces \u2013\n// \"Appendix\", \": \" and \"List of Data Prep/Cleaning tools\" \u2013 mirroring\n// the way the author typed it, and a \"_GoBack\" bookmark (Word's\n// auto-maintained \"last edit position\" marker) is left at the very end\n// of that text, matching where Word leaves it after the most recent\n// edit.\n\nconst body = context.document.body;\n\n// 1) \"Appendix: List of Data Prep/Cleaning tools\" heading, H1 style.\nconst heading = body.insertParagraph(\"Appendix\", Word.InsertLocation.end);\nheading.style = \"H1\";\nawait context.sync();\n\nheading.insertText(\": \", Word.InsertLocation.end);\nawait context.sync();\n\nheading.insertText(\"List of Data Prep/Cleaning tools\", Word.InsertLocation.end);\nawait context.sync();\n\n// Word leaves its \"_GoBack\" bookmark at the location of the most recent\n// edit -- i.e. the end of the text we just typed.\nconst endOfHeading = heading.getRange(Word.RangeLocation.end);\nendOfHeading.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) A trailing, empty paragraph using the \"Para\" style.\nconst trailer = body.insertParagraph(\"\", Word.InsertLocation.end);\ntrailer.style = \"Para\";\nawait context.sync();\n", "ps1": "# Commit: \"adding in appendix for data prep tools\"\n#\n# Appends a new \"Appendix\" heading (style H1) followed by an empty\n# paragraph (style \"Para\") to the very end of the document body, after\n# the last existing paragraph (\"...Reiterating that data analysis\n# assists our thinking, it does not replace it\").\n\n$d = $word.ActiveDocument\n\n# 1) \"Appendix: List of Data Prep/Cleaning tools\" heading, H1 style.\n$heading = $d.Content.Paragraphs.Add()\n$headingRange = $heading.Range\n$headingRange.Text = \"Appendix: List of Data Prep/Cleaning tools\"\n$headingRange.Style = \"H1\"\n\n# Word leaves its \"_GoBack\" bookmark at the location of the most recent\n# edit -- i.e. the end of the text we just typed.\n$d.Bookmarks.Add(\"_GoBack\", $headingRange)\n\n# 2) A trailing, empty paragraph using the \"Para\" style.\n$trailer = $d.Content.Paragraphs.Add()\n$trailer.Range.Style = \"Para\"\n"}
